$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet: Rev. C -> Rev. D (this also updates the suffixed
# _xlnm.Print_Area_0... defined names automatically)
$ws.Name = "Domino JTAG SPI Rev. D"

# Re-apply the print area so the plain _xlnm.Print_Area defined name
# also gets refreshed to point at the renamed sheet.
$ws.PageSetup.PrintArea = '$A$1:$I$12'

# Column width tweaks (silkscreen font ratio change caused Calc to
# recompute the autofit column widths slightly wider).
$ws.Range("A1").EntireColumn.ColumnWidth = 4.166666666666667
$ws.Range("B1").EntireColumn.ColumnWidth = 4.166666666666667
$ws.Range("C1").EntireColumn.ColumnWidth = 26.166666666666668
$ws.Range("D1").EntireColumn.ColumnWidth = 28.166666666666668
$ws.Range("E1").EntireColumn.ColumnWidth = 31.333333333333332
$ws.Range("F1").EntireColumn.ColumnWidth = 26.166666666666668
$ws.Range("G1").EntireColumn.ColumnWidth = 52.166666666666664
$ws.Range("H1").EntireColumn.ColumnWidth = 70.0
$ws.Range("I1").EntireColumn.ColumnWidth = 11.0
